# Weekly update: a new price-report entry (Primera / Segunda) is inserted
# at the top of the Cilantro data block (rows 163-164), pushing every
# subsequent entry down by one pair of rows. The new entry duplicates the
# row that used to be first in the block (now shifted to rows 165-166),
# except for its date (column D), which is the new reporting date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 163 - this pushes the existing data
# (formerly rows 163:202) down to rows 165:204.
$ws.Rows.Item(163).Resize(2).Insert()

# The row that used to be on top of the block is now at 165:166. Copy its
# values/formatting into the freshly inserted rows 163:164 so the new
# entry starts out identical to it.
$ws.Range("A165:R166").Copy()
$ws.Range("A163:R164").PasteSpecial(-4104)
$excel.CutCopyMode = $false

# Finally, set the new entry's date to the new reporting date.
$ws.Range("D163").Value = 44736
$ws.Range("D164").Value = 44736
